# Gestion des cotisations | enregistrement de paiement avec gestion des comptes de cotisations
#
# The "Beneficiaire" sheet holds two rows of people. This edit replaces the
# previously free-typed text values for "Date de naissance" / "Date d'entree"
# with real Excel date values, clears out the "N deg CIN" text value (the
# cell keeps its number format, it is simply emptied), and tweaks the
# "N deg Matricule" numbers for both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 already carries the date number-format (style index 4 in the original
# file) - use it as the format donor so the destination cells end up on the
# exact same style instead of Excel minting a brand new custom number
# format for an equivalent "m/d/yyyy" pattern.
$dateFormatDonor = $ws.Cells.Item(2, 5)
$dateFormatDonor.Copy()

# ----- Row 2 (Gabeta Sandratra) -----
$ws.Cells.Item(2, 1).Value = 6              # N deg Matricule: 4 -> 6
$ws.Cells.Item(2, 5).Value = 34742          # Date de naissance -> 12/02/1995
$ws.Cells.Item(2, 6).Value = ""             # N deg CIN cleared (format kept)
$ws.Cells.Item(2, 8).PasteSpecial(-4122)    # apply the date format (xlPasteFormats)
$ws.Cells.Item(2, 8).Value = 43263          # Date d'entree -> 12/06/2018

# ----- Row 3 (Rakoto Gabeta) -----
$ws.Cells.Item(3, 1).Value = 12             # N deg Matricule: 13 -> 12
$ws.Cells.Item(3, 5).PasteSpecial(-4122)    # apply the date format (xlPasteFormats)
$ws.Cells.Item(3, 5).Value = 37024          # Date de naissance -> 13/05/2001
$ws.Cells.Item(3, 6).Value = ""             # N deg CIN cleared (format kept)
$ws.Cells.Item(3, 8).PasteSpecial(-4122)    # apply the date format (xlPasteFormats)
$ws.Cells.Item(3, 8).Value = 43994          # Date d'entree -> 12/06/2020

# Matches the selection left behind in the saved file.
$ws.Range("G8").Select() | Out-Null
